$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "Símbulo" -> "Símbolo" in the header cell C1
$ws.Range("C1").Value = "Símbolo"

# Update the selected cell to C1
$ws.Range("C1").Select()
